$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B4 image path string value
$ws.Range("B4").Value = "Assets/Res/UI/Sprite/s1.png"

# Add new row of data
$ws.Range("A5").Value = "TextureA"
$ws.Range("B5").Value = "Assets/Res/UI/Texture/t1.png"

# Move the active selection to B6 (as if the user pressed Enter after editing B5)
$ws.Range("B6").Select()
